$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each row, copy the format from the Q column cell into the new R column
# cell, then set the R cell value (year for row 3, data value otherwise).
# Doing the format-copy before the value assignment lets the paste reuse the
# existing style index (matching Q) instead of minting a new one.
$ws.Range("Q3").Copy() | Out-Null
$ws.Range("R3").PasteSpecial(-4122) | Out-Null
$ws.Range("R3").Value = 2021
$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null
$ws.Range("R4").Value = 58.14349653559799
$ws.Range("Q5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null
$ws.Range("R5").Value = 50.405857641278807
$ws.Range("Q6").Copy() | Out-Null
$ws.Range("R6").PasteSpecial(-4122) | Out-Null
$ws.Range("R6").Value = 65.995789757646122
$ws.Range("Q7").Copy() | Out-Null
$ws.Range("R7").PasteSpecial(-4122) | Out-Null
$ws.Range("R7").Value = 47.339416388110941
$ws.Range("Q8").Copy() | Out-Null
$ws.Range("R8").PasteSpecial(-4122) | Out-Null
$ws.Range("R8").Value = 44.18457369250482
$ws.Range("Q9").Copy() | Out-Null
$ws.Range("R9").PasteSpecial(-4122) | Out-Null
$ws.Range("R9").Value = 50.379263611270765
$ws.Range("Q10").Copy() | Out-Null
$ws.Range("R10").PasteSpecial(-4122) | Out-Null
$ws.Range("R10").Value = 54.819947539591084
$ws.Range("Q11").Copy() | Out-Null
$ws.Range("R11").PasteSpecial(-4122) | Out-Null
$ws.Range("R11").Value = 47.679920417302263
$ws.Range("Q12").Copy() | Out-Null
$ws.Range("R12").PasteSpecial(-4122) | Out-Null
$ws.Range("R12").Value = 61.861274529713718
$ws.Range("Q13").Copy() | Out-Null
$ws.Range("R13").PasteSpecial(-4122) | Out-Null
$ws.Range("R13").Value = 36.712395096811576
$ws.Range("Q14").Copy() | Out-Null
$ws.Range("R14").PasteSpecial(-4122) | Out-Null
$ws.Range("R14").Value = 26.872053459579295
$ws.Range("Q15").Copy() | Out-Null
$ws.Range("R15").PasteSpecial(-4122) | Out-Null
$ws.Range("R15").Value = 46.638444428499682
$ws.Range("Q16").Copy() | Out-Null
$ws.Range("R16").PasteSpecial(-4122) | Out-Null
$ws.Range("R16").Value = 51.155081745820631
$ws.Range("Q17").Copy() | Out-Null
$ws.Range("R17").PasteSpecial(-4122) | Out-Null
$ws.Range("R17").Value = 43.08338023862634
$ws.Range("Q18").Copy() | Out-Null
$ws.Range("R18").PasteSpecial(-4122) | Out-Null
$ws.Range("R18").Value = 58.934228062068456
$ws.Range("Q19").Copy() | Out-Null
$ws.Range("R19").PasteSpecial(-4122) | Out-Null
$ws.Range("R19").Value = 54.51979816984521
$ws.Range("Q20").Copy() | Out-Null
$ws.Range("R20").PasteSpecial(-4122) | Out-Null
$ws.Range("R20").Value = 52.474443936678909
$ws.Range("Q21").Copy() | Out-Null
$ws.Range("R21").PasteSpecial(-4122) | Out-Null
$ws.Range("R21").Value = 56.519551395440942
$ws.Range("Q22").Copy() | Out-Null
$ws.Range("R22").PasteSpecial(-4122) | Out-Null
$ws.Range("R22").Value = 46.970408642555192
$ws.Range("Q23").Copy() | Out-Null
$ws.Range("R23").PasteSpecial(-4122) | Out-Null
$ws.Range("R23").Value = 27.43769048802011
$ws.Range("Q24").Copy() | Out-Null
$ws.Range("R24").PasteSpecial(-4122) | Out-Null
$ws.Range("R24").Value = 66.104415920267911
$ws.Range("Q25").Copy() | Out-Null
$ws.Range("R25").PasteSpecial(-4122) | Out-Null
$ws.Range("R25").Value = 88.246666265390886
$ws.Range("Q26").Copy() | Out-Null
$ws.Range("R26").PasteSpecial(-4122) | Out-Null
$ws.Range("R26").Value = 71.914698721605745
$ws.Range("Q27").Copy() | Out-Null
$ws.Range("R27").PasteSpecial(-4122) | Out-Null
$ws.Range("R27").Value = 105.10059183863845
$ws.Range("Q28").Copy() | Out-Null
$ws.Range("R28").PasteSpecial(-4122) | Out-Null
$ws.Range("R28").Value = 63.980940123966526
$ws.Range("Q29").Copy() | Out-Null
$ws.Range("R29").PasteSpecial(-4122) | Out-Null
$ws.Range("R29").Value = 55.546587096180644
$ws.Range("Q30").Copy() | Out-Null
$ws.Range("R30").PasteSpecial(-4122) | Out-Null
$ws.Range("R30").Value = 73.505198287622903
$ws.Range("Q31").Copy() | Out-Null
$ws.Range("R31").PasteSpecial(-4122) | Out-Null
$ws.Range("R31").Value = 43.916363725083563
$ws.Range("Q32").Copy() | Out-Null
$ws.Range("R32").PasteSpecial(-4122) | Out-Null
$ws.Range("R32").Value = 40.980198843051781
$ws.Range("Q33").Copy() | Out-Null
$ws.Range("R33").PasteSpecial(-4122) | Out-Null
$ws.Range("R33").Value = 47.015458682814909

$excel.CutCopyMode = 0

# Update the active selection to match the authored state.
$ws.Range("T3").Select() | Out-Null

Write-Output "done"
